$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 108.595075
$ws.Range("H2").Value = 325.785225
$ws.Range("I2").Value = 0.04639022893696803
$ws.Range("J2").Value = 0.04639022893696803
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 883.0119660725499
$ws.Range("R2").Value = 7947.107694652948
$ws.Range("S2").Value = 0.0009700176333839064
$ws.Range("T2").Value = 0.0009700176333839066

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 108.595075
$ws.Range("H3").Value = 325.785225
$ws.Range("I3").Value = 0.04639022893696803
$ws.Range("J3").Value = 0.04639022893696803
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 26429.47826643657
$ws.Range("R3").Value = 237865.3043979291
$ws.Range("S3").Value = 0.0290336495366064
$ws.Range("T3").Value = 0.0290336495366064

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 108.595075
$ws.Range("H4").Value = 325.785225
$ws.Range("I4").Value = 0.04639022893696803
$ws.Range("J4").Value = 0.04639022893696803
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 11287.65860550621
$ws.Range("R4").Value = 101588.9274495559
$ws.Range("S4").Value = 0.01239986354393187
$ws.Range("T4").Value = 0.01239986354393187

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 108.595075
$ws.Range("H5").Value = 325.785225
$ws.Range("I5").Value = 0.04639022893696803
$ws.Range("J5").Value = 0.04639022893696803
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 3629.111590259533
$ws.Range("R5").Value = 32662.0043123358
$ws.Range("S5").Value = 0.00398669822304585
$ws.Range("T5").Value = 0.003986698223045851

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 251.7279513333333
$ws.Range("H6").Value = 755.183854
$ws.Range("I6").Value = 0.1075345018380187
$ws.Range("J6").Value = 0.1075345018380187
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 2046.858876631945
$ws.Range("R6").Value = 18421.72988968751
$ws.Range("S6").Value = 0.002248541672897589
$ws.Range("T6").Value = 0.00224854167289759

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 251.7279513333333
$ws.Range("H7").Value = 755.183854
$ws.Range("I7").Value = 0.1075345018380187
$ws.Range("J7").Value = 0.1075345018380187
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 61264.6422392446
$ws.Range("R7").Value = 551381.7801532014
$ws.Range("S7").Value = 0.06730122077432987
$ws.Range("T7").Value = 0.0673012207743299

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 251.7279513333333
$ws.Range("H8").Value = 755.183854
$ws.Range("I8").Value = 0.1075345018380187
$ws.Range("J8").Value = 0.1075345018380187
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 26165.2674038316
$ws.Range("R8").Value = 235487.4066344844
$ws.Range("S8").Value = 0.02874340523018062
$ws.Range("T8").Value = 0.02874340523018063

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 251.7279513333333
$ws.Range("H9").Value = 755.183854
$ws.Range("I9").Value = 0.1075345018380187
$ws.Range("J9").Value = 0.1075345018380187
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 8412.433305802198
$ws.Range("R9").Value = 75711.89975221979
$ws.Range("S9").Value = 0.009241334160610621
$ws.Range("T9").Value = 0.009241334160610624

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1894.44458
$ws.Range("H10").Value = 5683.33374
$ws.Range("I10").Value = 0.8092790375125046
$ws.Range("J10").Value = 0.8092790375125047
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 15404.17218001172
$ws.Range("R10").Value = 138637.5496201055
$ws.Range("S10").Value = 0.01692198884773153
$ws.Range("T10").Value = 0.01692198884773153

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1894.44458
$ws.Range("H11").Value = 5683.33374
$ws.Range("I11").Value = 0.8092790375125046
$ws.Range("J11").Value = 0.8092790375125047
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 461063.099353986
$ws.Range("R11").Value = 4149567.894185874
$ws.Range("S11").Value = 0.5064929510131422
$ws.Range("T11").Value = 0.5064929510131424

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1894.44458
$ws.Range("H12").Value = 5683.33374
$ws.Range("I12").Value = 0.8092790375125046
$ws.Range("J12").Value = 0.8092790375125047
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 196913.5678214836
$ws.Range("R12").Value = 1772222.110393352
$ws.Range("S12").Value = 0.2163160187839212
$ws.Range("T12").Value = 0.2163160187839213

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1894.44458
$ws.Range("H13").Value = 5683.33374
$ws.Range("I13").Value = 0.8092790375125046
$ws.Range("J13").Value = 0.8092790375125047
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 63309.96854491194
$ws.Range("R13").Value = 569789.7169042075
$ws.Range("S13").Value = 0.06954807886770964
$ws.Range("T13").Value = 0.06954807886770967

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 86.13644799999999
$ws.Range("H14").Value = 258.409344
$ws.Range("I14").Value = 0.03679623171250852
$ws.Range("J14").Value = 0.03679623171250853
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 700.3956146168318
$ws.Range("R14").Value = 6303.560531551487
$ws.Range("S14").Value = 0.0007694075761451974
$ws.Range("T14").Value = 0.0007694075761451976

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 86.13644799999999
$ws.Range("H15").Value = 258.409344
$ws.Range("I15").Value = 0.03679623171250852
$ws.Range("J15").Value = 0.03679623171250853
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 20963.5785081786
$ws.Range("R15").Value = 188672.2065736074
$ws.Range("S15").Value = 0.02302917921056844
$ws.Range("T15").Value = 0.02302917921056845

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 86.13644799999999
$ws.Range("H16").Value = 258.409344
$ws.Range("I16").Value = 0.03679623171250852
$ws.Range("J16").Value = 0.03679623171250853
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 8953.249661781972
$ws.Range("R16").Value = 80579.24695603775
$ws.Range("S16").Value = 0.009835438682269737
$ws.Range("T16").Value = 0.009835438682269741

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 86.13644799999999
$ws.Range("H17").Value = 258.409344
$ws.Range("I17").Value = 0.03679623171250852
$ws.Range("J17").Value = 0.03679623171250853
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 2878.572364175701
$ws.Range("R17").Value = 25907.15127758131
$ws.Range("S17").Value = 0.003162206243525144
$ws.Range("T17").Value = 0.003162206243525145

Write-Host "Updated cells"